$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.407.46"
$ws.Range("E2").Value = "  +5.82%  "
$ws.Range("D3").Value = "2.040.80"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.63"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("E6").Value = "  +3.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.10"
$ws.Range("E7").Value = "  +18.25%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +6.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.09"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0754"
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.907"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.14"
$ws.Range("E14").Value = "  +8.15%  "
$ws.Range("D15").Value = "2.338.75"
$ws.Range("E15").Value = "  +3.08%  "
$ws.Range("E16").Value = "  +9.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.70"
$ws.Range("E17").Value = "  +23.68%  "
$ws.Range("D18").Value = "2.042.04"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("D19").Value = "37.287.84"
$ws.Range("E19").Value = "  +5.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.14"
$ws.Range("E20").Value = "  +5.53%  "
$ws.Range("D21").Value = "0.0₃0873"
$ws.Range("E21").Value = "  +5.82%  "
$ws.Range("E22").Value = "  +8.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.82"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("E24").Value = "  +23.39%  "
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.57"
$ws.Range("E27").Value = "  +6.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.90"
$ws.Range("E28").Value = "  +2.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.88"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.24"
$ws.Range("E30").Value = "  +10.89%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.23"
$ws.Range("E32").Value = "  +10.23%  "
$ws.Range("E33").Value = "  +27.67%  "
$ws.Range("E34").Value = "  +12.99%  "
$ws.Range("E35").Value = "  +6.26%  "
$ws.Range("E36").Value = "  +13.93%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.96"
$ws.Range("E39").Value = "  +24.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  +19.45%  "
$ws.Range("E41").Value = "  +5.57%  "
$ws.Range("E42").Value = "  +4.05%  "
$ws.Range("E43").Value = "  +6.50%  "
$ws.Range("E44").Value = "  +22.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.17"
$ws.Range("E45").Value = "  +12.02%  "
$ws.Range("E46").Value = "  +7.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.14"
$ws.Range("E47").Value = "  +12.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.42"
$ws.Range("E48").Value = "  +7.70%  "
$ws.Range("D49").Value = "1.421.45"
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.57"
$ws.Range("E51").Value = "  +7.44%  "
